$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '30.323.68'
$ws.Cells.Item(2, 5).Value = '  +0.06%  '
$ws.Cells.Item(3, 4).Value = '1.842.87'
$ws.Cells.Item(3, 5).Value = '  -0.84%  '
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '0.9996'
$ws.Cells.Item(4, 5).Value = '  -0.35%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '233.11'
$ws.Cells.Item(5, 5).Value = '  +0.31%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.9992'
$ws.Cells.Item(6, 5).Value = '  -0.29%  '
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.4659'
$ws.Cells.Item(7, 5).Value = '  -1.84%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.2731'
$ws.Cells.Item(8, 5).Value = '  -0.77%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.06286'
$ws.Cells.Item(9, 5).Value = '  -2.38%  '
$ws.Cells.Item(10, 4).Value = '1.836.37'
$ws.Cells.Item(10, 5).Value = '  -1.18%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.07418'
$ws.Cells.Item(11, 5).Value = '  +0.11%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '16.26'
$ws.Cells.Item(12, 5).Value = '  +1.48%  '
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '4.940'
$ws.Cells.Item(13, 5).Value = '  -1.09%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '83.76'
$ws.Cells.Item(14, 5).Value = '  -2.41%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.6206'
$ws.Cells.Item(15, 5).Value = '  -1.59%  '
$ws.Cells.Item(16, 4).Value = '30.264.53'
$ws.Cells.Item(16, 5).Value = '  -0.17%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '0.9987'
$ws.Cells.Item(17, 5).Value = '  -0.25%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '227.93'
$ws.Cells.Item(18, 5).Value = '  +0.38%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '0.000007299'
$ws.Cells.Item(19, 5).Value = '  -0.16%  '
$ws.Cells.Item(20, 5).Value = '  -3.80%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '0.9990'
$ws.Cells.Item(21, 5).Value = '  -0.55%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '4.914'
$ws.Cells.Item(22, 5).Value = '  -3.39%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '5.860'
$ws.Cells.Item(23, 5).Value = '  -3.15%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '9.181'
$ws.Cells.Item(24, 5).Value = '  -0.39%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '164.34'
$ws.Cells.Item(25, 5).Value = '  -1.74%  '
$ws.Cells.Item(26, 5).Value = '  -0.12%  '
$ws.Cells.Item(27, 5).Value = '  +0.07%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '0.1030'
$ws.Cells.Item(28, 5).Value = '  -0.45%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '1.370'
$ws.Cells.Item(29, 5).Value = '  -0.87%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '4.069'
$ws.Cells.Item(30, 5).Value = '  -3.24%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '3.804'
$ws.Cells.Item(31, 5).Value = '  -2.30%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '0.04839'
$ws.Cells.Item(32, 5).Value = '  -1.46%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '1.141'
$ws.Cells.Item(33, 5).Value = '  -0.89%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '0.7119'
$ws.Cells.Item(34, 5).Value = '  -1.72%  '
$ws.Cells.Item(35, 5).Value = '  -0.73%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.01900'
$ws.Cells.Item(36, 5).Value = '  +0.39%  '
$ws.Cells.Item(37, 5).Value = '  +1.17%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.8847'
$ws.Cells.Item(38, 5).Value = '  -2.28%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '105.23'
$ws.Cells.Item(39, 5).Value = '  +0.07%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '1.926'
$ws.Cells.Item(40, 5).Value = '  -2.67%  '
$ws.Cells.Item(41, 5).Value = '  +0.58%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '5.553'
$ws.Cells.Item(42, 5).Value = '  +0.32%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.4026'
$ws.Cells.Item(43, 5).Value = '  -1.71%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '7.158'
$ws.Cells.Item(44, 5).Value = '  +1.92%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '62.10'
$ws.Cells.Item(45, 5).Value = '  +2.01%  '
$ws.Cells.Item(46, 5).Value = '  -0.39%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '8.633'
$ws.Cells.Item(47, 5).Value = '  -1.48%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '33.21'
$ws.Cells.Item(48, 5).Value = '  +1.01%  '
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.05508'
$ws.Cells.Item(49, 5).Value = '  -2.01%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '1.350'
$ws.Cells.Item(50, 5).Value = '  -3.32%  '
$ws.Cells.Item(51, 5).Value = '  -1.72%  '
